$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before the existing F (recommended_level), shifting it to K
$ws.Range("F1:J1").EntireColumn.Insert()

# Header row
$ws.Range("F1").Value = "frequency"
$ws.Range("G1").Value = "frequency_occurrence"
$ws.Range("H1").Value = "frequency_occurrence_probab"
$ws.Range("I1").Value = "max_probab"
$ws.Range("J1").Value = "max_probab_percentage"

# Data rows: frequency, frequency_occurrence, frequency_occurrence_probab, max_probab, max_probab_percentage
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = '{"L3":1}'
$ws.Range("H2").Value = '{"L3":1.0}'
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "'100.00"
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = '{"L3":7,"L2":1}'
$ws.Range("H3").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I3").Value = 0.875
$ws.Range("J3").Value = "'87.50"
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = '{"L3":8}'
$ws.Range("H4").Value = '{"L3":1.0}'
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "'100.00"
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = '{"L3":8}'
$ws.Range("H5").Value = '{"L3":1.0}'
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = "'100.00"
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = '{"L2":7,"L3":1}'
$ws.Range("H6").Value = '{"L2":0.875,"L3":0.125}'
$ws.Range("I6").Value = 0.875
$ws.Range("J6").Value = "'87.50"
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = '{"L3":8}'
$ws.Range("H7").Value = '{"L3":1.0}'
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = "'100.00"
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = '{"L2":5,"L1":3}'
$ws.Range("H8").Value = '{"L2":0.625,"L1":0.375}'
$ws.Range("I8").Value = 0.625
$ws.Range("J8").Value = "'62.50"
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = '{"L2":8}'
$ws.Range("H9").Value = '{"L2":1.0}'
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = "'100.00"
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = '{"L3":8}'
$ws.Range("H10").Value = '{"L3":1.0}'
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "'100.00"
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = '{"L3":6,"L2":2}'
$ws.Range("H11").Value = '{"L3":0.75,"L2":0.25}'
$ws.Range("I11").Value = 0.75
$ws.Range("J11").Value = "'75.00"
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = '{"L2":6,"L1":2}'
$ws.Range("H12").Value = '{"L2":0.75,"L1":0.25}'
$ws.Range("I12").Value = 0.75
$ws.Range("J12").Value = "'75.00"
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = '{"L3":8}'
$ws.Range("H13").Value = '{"L3":1.0}'
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = "'100.00"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = '{"L3":6,"L2":2}'
$ws.Range("H14").Value = '{"L3":0.75,"L2":0.25}'
$ws.Range("I14").Value = 0.75
$ws.Range("J14").Value = "'75.00"
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = '{"L3":7,"L2":1}'
$ws.Range("H15").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I15").Value = 0.875
$ws.Range("J15").Value = "'87.50"
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = '{"L3":8}'
$ws.Range("H16").Value = '{"L3":1.0}'
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = "'100.00"
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = '{"L3":8}'
$ws.Range("H17").Value = '{"L3":1.0}'
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = "'100.00"
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = '{"L2":1}'
$ws.Range("H18").Value = '{"L2":1.0}'
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = "'100.00"
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = '{"L2":6,"L3":2}'
$ws.Range("H19").Value = '{"L2":0.75,"L3":0.25}'
$ws.Range("I19").Value = 0.75
$ws.Range("J19").Value = "'75.00"
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = '{"L2":6,"L3":2}'
$ws.Range("H20").Value = '{"L2":0.75,"L3":0.25}'
$ws.Range("I20").Value = 0.75
$ws.Range("J20").Value = "'75.00"
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = '{"L2":7,"L1":1}'
$ws.Range("H22").Value = '{"L2":0.875,"L1":0.125}'
$ws.Range("I22").Value = 0.875
$ws.Range("J22").Value = "'87.50"
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = '{"L3":1}'
$ws.Range("H23").Value = '{"L3":1.0}'
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = "'100.00"
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = '{"L3":1}'
$ws.Range("H24").Value = '{"L3":1.0}'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = "'100.00"
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = '{"L2":7,"L3":1}'
$ws.Range("H25").Value = '{"L2":0.875,"L3":0.125}'
$ws.Range("I25").Value = 0.875
$ws.Range("J25").Value = "'87.50"
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = '{"L3":8}'
$ws.Range("H26").Value = '{"L3":1.0}'
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = "'100.00"
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = '{"L3":6,"L2":2}'
$ws.Range("H27").Value = '{"L3":0.75,"L2":0.25}'
$ws.Range("I27").Value = 0.75
$ws.Range("J27").Value = "'75.00"
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = '{"L3":7,"L2":1}'
$ws.Range("H28").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I28").Value = 0.875
$ws.Range("J28").Value = "'87.50"
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = '{"L2":1}'
$ws.Range("H35").Value = '{"L2":1.0}'
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = "'100.00"
$ws.Range("F38").Value = 8
$ws.Range("G38").Value = '{"L1":8}'
$ws.Range("H38").Value = '{"L1":1.0}'
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = "'100.00"
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = '{"L1":1}'
$ws.Range("H39").Value = '{"L1":1.0}'
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = "'100.00"
